# Auto-generated Excel COM-interop script
# Rewrites FBS (sheet1) and Other (sheet2) tables to the refreshed
# cfb_weather data snapshot (Timestamp 2025-11-30T16:22:53.319472).

$wb = $excel.ActiveWorkbook

# ---------------- FBS sheet ----------------
$wsFBS = $wb.Worksheets.Item("FBS")

$data1 = New-Object 'object[,]' 7,37
$data1[0,0] = 'Game'
$data1[0,1] = 'Date'
$data1[0,2] = 'Time'
$data1[0,3] = 'wind_vol'
$data1[0,4] = 'orient'
$data1[0,5] = 'wind_impact'
$data1[0,6] = 'weakest_wind_effect'
$data1[0,7] = 'travel_alt'
$data1[0,8] = 'home_temp'
$data1[0,9] = 'away_temp'
$data1[0,10] = 'wind_avg'
$data1[0,11] = 'year_built'
$data1[0,12] = 'wind_dir_1h'
$data1[0,13] = 'wind_dir_2h'
$data1[0,14] = 'temp_fg'
$data1[0,15] = 'wind_fg'
$data1[0,16] = 'wind_dir_fg'
$data1[0,17] = 'rain_fg'
$data1[0,18] = 'gs_fg'
$data1[0,19] = 'away_fg'
$data1[0,20] = 'wind_diff'
$data1[0,21] = 'game_loc'
$data1[0,22] = 'Fd_open'
$data1[0,23] = 'Odds_o'
$data1[0,24] = 'FD_now'
$data1[0,25] = 'Odds_n'
$data1[0,26] = 'Open'
$data1[0,27] = 'Current'
$data1[0,28] = 'Spread'
$data1[0,29] = 'Total_proj'
$data1[0,30] = 'Move_t'
$data1[0,31] = 'Move_s'
$data1[0,32] = 'My_total'
$data1[0,33] = 'Edge'
$data1[0,34] = 'My_spread'
$data1[0,35] = 'Edge_s'
$data1[0,36] = 'Timestamp'
$data1[1,0] = 'UNLV @ Boise State'
$data1[1,1] = 'FRI 12/05'
$data1[1,2] = '06:00 PM'
$data1[1,3] = 'High'
$data1[1,4] = 'N-S'
$data1[1,5] = 'Med'
$data1[1,6] = 'E'
$data1[1,8] = 53.65
$data1[1,9] = 70.04000000000001
$data1[1,10] = 6.8
$data1[1,11] = 1970
$data1[1,12] = 'ESE'
$data1[1,13] = 'ESE'
$data1[1,14] = 48.5
$data1[1,15] = 18.2
$data1[1,16] = 'ESE'
$data1[1,17] = 0
$data1[1,18] = -6.5
$data1[1,19] = 0
$data1[1,20] = 11.4
$data1[1,21] = '43.6028839, -116.1958882'
$data1[1,22] = 56.5
$data1[1,23] = -110
$data1[1,24] = 57.5
$data1[1,25] = -110
$data1[1,26] = -3
$data1[1,27] = -3.5
$data1[1,30] = 0.01769911504424779
$data1[1,31] = 0.5
$data1[1,36] = '2025-11-30T16:22:53.319472'
$data1[2,0] = 'Troy @ James Madison'
$data1[2,1] = 'FRI 12/05'
$data1[2,2] = '07:00 PM'
$data1[2,3] = 'High'
$data1[2,4] = 'N-S'
$data1[2,5] = 'Med'
$data1[2,6] = 'E'
$data1[2,7] = 244.5283813
$data1[2,8] = 54.52
$data1[2,9] = 65.95999999999999
$data1[2,10] = 5.9
$data1[2,11] = 1975
$data1[2,12] = 'N'
$data1[2,13] = 'N'
$data1[2,14] = 24.32
$data1[2,15] = 7
$data1[2,16] = 'N'
$data1[2,17] = 0.1
$data1[2,18] = -0.71
$data1[2,19] = -0.96
$data1[2,20] = 1.1
$data1[2,21] = '38.4352919, -78.8729349'
$data1[2,22] = 47.5
$data1[2,23] = -105
$data1[2,24] = 47.5
$data1[2,25] = -105
$data1[2,26] = -21
$data1[2,27] = -22
$data1[2,30] = 0
$data1[2,31] = 1
$data1[2,36] = '2025-11-30T16:22:53.319472'
$data1[3,0] = 'Kennesaw State @ Jacksonville State'
$data1[3,1] = 'FRI 12/05'
$data1[3,2] = '06:00 PM'
$data1[3,3] = 'Low'
$data1[3,4] = 'E-W'
$data1[3,5] = 'High'
$data1[3,6] = 'N'
$data1[3,7] = -98.89108280000002
$data1[3,8] = 63.15
$data1[3,9] = 61.32
$data1[3,10] = 4.7
$data1[3,11] = 1947
$data1[3,12] = 'W'
$data1[3,13] = 'W'
$data1[3,14] = 43.58000000000001
$data1[3,15] = 5.8
$data1[3,16] = 'W'
$data1[3,17] = 0.3
$data1[3,18] = 0
$data1[3,19] = 0
$data1[3,20] = 1.1
$data1[3,21] = '33.8201052, -85.76647'
$data1[3,22] = 58.5
$data1[3,23] = -110
$data1[3,24] = 58.5
$data1[3,25] = -105
$data1[3,26] = -1.5
$data1[3,27] = -1.5
$data1[3,30] = 0
$data1[3,31] = 0
$data1[3,36] = '2025-11-30T16:22:53.319472'
$data1[4,0] = 'North Texas @ Tulane'
$data1[4,1] = 'FRI 12/05'
$data1[4,2] = '07:00 PM'
$data1[4,3] = 'High'
$data1[4,4] = 'NE-SW'
$data1[4,5] = 'Med'
$data1[4,6] = 'NW'
$data1[4,7] = -198.6871948
$data1[4,8] = 70.11
$data1[4,9] = 66.45
$data1[4,10] = 10.1
$data1[4,11] = 2014
$data1[4,12] = 'SSE'
$data1[4,13] = 'SSE'
$data1[4,14] = 60.2
$data1[4,15] = 6.7
$data1[4,16] = 'SSE'
$data1[4,17] = 0
$data1[4,18] = 0
$data1[4,19] = 0
$data1[4,20] = -3.4
$data1[4,21] = '29.944616, -90.116692'
$data1[4,22] = 62.5
$data1[4,23] = -110
$data1[4,24] = 67.5
$data1[4,25] = -115
$data1[4,26] = 3
$data1[4,27] = 2.5
$data1[4,30] = 0.08
$data1[4,31] = 0.5
$data1[4,36] = '2025-11-30T16:22:53.319472'
$data1[5,0] = 'Miami (OH) @ Western Michigan'
$data1[5,1] = 'SAT 12/06'
$data1[5,2] = '12:00 PM'
$data1[5,3] = 'Mid'
$data1[5,4] = 'NE-SW'
$data1[5,5] = 'High'
$data1[5,6] = 'SW'
$data1[5,7] = 1.556762700000007
$data1[5,8] = 50.43
$data1[5,9] = 54.38
$data1[5,10] = 11.7
$data1[5,11] = 1939
$data1[5,12] = 'NNE'
$data1[5,13] = 'NNE'
$data1[5,14] = 30.5
$data1[5,15] = 6.5
$data1[5,16] = 'NNE'
$data1[5,17] = 0
$data1[5,18] = 0
$data1[5,19] = 0
$data1[5,20] = -5.2
$data1[5,21] = '42.2860064, -85.6007573'
$data1[5,22] = 42.5
$data1[5,23] = -110
$data1[5,24] = 43.5
$data1[5,25] = -105
$data1[5,30] = 0.02352941176470588
$data1[5,36] = '2025-11-30T16:22:53.319472'
$data1[6,0] = 'Duke @ Virginia'
$data1[6,1] = 'SAT 12/06'
$data1[6,2] = '08:00 PM'
$data1[6,3] = 'High'
$data1[6,4] = 'NW-SE'
$data1[6,5] = 'Med'
$data1[6,7] = 67.46492769999999
$data1[6,8] = 58.23
$data1[6,9] = 61.08
$data1[6,10] = 4.5
$data1[6,11] = 1931
$data1[6,12] = 'SSW'
$data1[6,13] = 'SSE'
$data1[6,14] = 32.12
$data1[6,15] = 1.2
$data1[6,16] = 'SSE'
$data1[6,17] = 0
$data1[6,18] = 0
$data1[6,19] = 0
$data1[6,20] = -3.3
$data1[6,21] = '38.0311801, -78.5137897'
$data1[6,22] = 58.5
$data1[6,23] = -110
$data1[6,24] = 58.5
$data1[6,25] = -115
$data1[6,26] = -2.5
$data1[6,27] = -3
$data1[6,30] = 0
$data1[6,31] = 0.5
$data1[6,36] = '2025-11-30T16:22:53.319472'

$wsFBS.Range("A1:AK7").Value = $data1

# ---------------- Other sheet ----------------
$wsOther = $wb.Worksheets.Item("Other")

$data2 = New-Object 'object[,]' 8,24
$data2[0,0] = 'Game'
$data2[0,1] = 'Home Team'
$data2[0,2] = 'Away Team'
$data2[0,3] = 'Date'
$data2[0,4] = 'Time'
$data2[0,5] = 'wind_vol'
$data2[0,6] = 'orient'
$data2[0,7] = 'wind_impact'
$data2[0,8] = 'weakest_wind_effect'
$data2[0,9] = 'travel_alt'
$data2[0,10] = 'home_temp'
$data2[0,11] = 'away_temp'
$data2[0,12] = 'wind_avg'
$data2[0,13] = 'year_built'
$data2[0,14] = 'wind_dir_1h'
$data2[0,15] = 'wind_dir_2h'
$data2[0,16] = 'temp_fg'
$data2[0,17] = 'wind_fg'
$data2[0,18] = 'wind_dir_fg'
$data2[0,19] = 'rain_fg'
$data2[0,20] = 'gs_fg'
$data2[0,21] = 'away_fg'
$data2[0,22] = 'wind_diff'
$data2[0,23] = 'game_loc'
$data2[1,0] = 'North Dakota vs Tarleton State'
$data2[1,1] = 'Tarleton State'
$data2[1,2] = 'North Dakota'
$data2[1,3] = 'SAT 12/06'
$data2[1,4] = '12:00 PM'
$data2[1,5] = 'Low'
$data2[1,10] = 66.93000000000001
$data2[1,11] = 42.5
$data2[1,14] = 'NE'
$data2[1,15] = 'NE'
$data2[1,16] = 68.54000000000001
$data2[1,17] = 25.5
$data2[1,18] = 'NE'
$data2[1,19] = 0
$data2[1,20] = -10
$data2[1,21] = 0
$data2[1,23] = '32.2191836, -98.2130634'
$data2[2,0] = 'Villanova vs Lehigh'
$data2[2,1] = 'Lehigh'
$data2[2,2] = 'Villanova'
$data2[2,3] = 'SAT 12/06'
$data2[2,4] = '12:00 PM'
$data2[2,5] = 'High'
$data2[2,9] = -37.06062315000001
$data2[2,10] = 54.29
$data2[2,11] = 55.05
$data2[2,13] = 1988
$data2[2,14] = 'SSE'
$data2[2,15] = 'SE'
$data2[2,16] = 30.44
$data2[2,17] = 1.9
$data2[2,18] = 'SE'
$data2[2,19] = 1.3
$data2[2,20] = -1.5
$data2[2,21] = 0
$data2[2,23] = '40.5890837, -75.3553874'
$data2[3,0] = 'Yale vs Montana State'
$data2[3,1] = 'Montana State'
$data2[3,2] = 'Yale'
$data2[3,3] = 'SAT 12/06'
$data2[3,4] = '12:00 PM'
$data2[3,5] = 'High'
$data2[3,9] = 1502.206045159
$data2[3,10] = 42.68
$data2[3,11] = 53.64
$data2[3,13] = 1973
$data2[3,14] = 'ESE'
$data2[3,15] = 'ESE'
$data2[3,16] = 21.5
$data2[3,17] = 2.2
$data2[3,18] = 'ESE'
$data2[3,19] = 0
$data2[3,20] = -1.06
$data2[3,21] = -3.5
$data2[3,23] = '45.659048, -111.049547'
$data2[4,0] = 'South Dakota vs Mercer'
$data2[4,1] = 'Mercer'
$data2[4,2] = 'South Dakota'
$data2[4,3] = 'SAT 12/06'
$data2[4,4] = '12:00 PM'
$data2[4,5] = 'Low'
$data2[4,9] = -234.7229156
$data2[4,10] = 64.83
$data2[4,11] = 51.08
$data2[4,13] = 2013
$data2[4,14] = 'SE'
$data2[4,15] = 'ESE'
$data2[4,16] = 54.26
$data2[4,17] = 3
$data2[4,18] = 'ESE'
$data2[4,19] = 0
$data2[4,20] = 0
$data2[4,21] = 0
$data2[4,23] = '32.8262075, -83.6522485'
$data2[5,0] = 'Abilene Christian vs Stephen F. Austin'
$data2[5,1] = 'Stephen F. Austin'
$data2[5,2] = 'Abilene Christian'
$data2[5,3] = 'SAT 12/06'
$data2[5,4] = '12:00 PM'
$data2[5,5] = 'Low'
$data2[5,9] = -430.90566101
$data2[5,10] = 68.06999999999999
$data2[5,11] = 67.58
$data2[5,13] = 1973
$data2[5,14] = 'NNE'
$data2[5,15] = 'NNE'
$data2[5,16] = 69.25999999999999
$data2[5,17] = 10.5
$data2[5,18] = 'NNE'
$data2[5,19] = 0
$data2[5,20] = 0
$data2[5,21] = 0
$data2[5,23] = '31.625719, -94.6444034'
$data2[6,0] = 'South Dakota State vs Montana'
$data2[6,1] = 'Montana'
$data2[6,2] = 'South Dakota State'
$data2[6,3] = 'SAT 12/06'
$data2[6,4] = '12:00 PM'
$data2[6,5] = 'High'
$data2[6,9] = 474.5684815
$data2[6,10] = 47.64
$data2[6,11] = 46.7
$data2[6,13] = 1986
$data2[6,14] = 'NE'
$data2[6,15] = 'ENE'
$data2[6,16] = 37.64
$data2[6,17] = 6.8
$data2[6,18] = 'ENE'
$data2[6,19] = 0.3
$data2[6,20] = 0
$data2[6,21] = 0
$data2[6,23] = '46.8638753, -113.9815042'
$data2[7,0] = 'Rhode Island vs UC Davis'
$data2[7,1] = 'UC Davis'
$data2[7,2] = 'Rhode Island'
$data2[7,3] = 'SAT 12/06'
$data2[7,4] = '12:00 PM'
$data2[7,5] = 'High'
$data2[7,9] = -21.30542278
$data2[7,10] = 62.21
$data2[7,11] = 52.81
$data2[7,13] = 2007
$data2[7,14] = 'NNE'
$data2[7,15] = 'NNE'
$data2[7,16] = 52.34
$data2[7,17] = 10.6
$data2[7,18] = 'NNE'
$data2[7,19] = 0
$data2[7,20] = 0
$data2[7,21] = 0
$data2[7,23] = '38.5365266, -121.7627936'

$wsOther.Range("A1:X8").Value = $data2

Write-Output "Updated FBS (7 rows) and Other (8 rows) with refreshed weather/odds data."
